$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# The edit appends one more bullet to the trailing numbered list:
#   "When listing empty players list exception fixed. Add Game constructor.
#    Renderer.cs left with default constructor (used parameter injection
#    instead) "
# with "Renderer.cs" flagged as a spell-check exception (proofErr), and the
# "_GoBack" bookmark relocated from the end of the previous ("... instead.")
# paragraph into the middle of the new one.
# ---------------------------------------------------------------------------

# 1. Detach the existing "_GoBack" bookmark from its current home (the end
#    of the last paragraph, "... Renderer.cs instead."). It will be
#    reinserted inside the brand-new paragraph below.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Append a new paragraph after the current last one. Word automatically
#    carries over the same paragraph style/list numbering (pStyle "a3",
#    numPr ilvl 0 / numId 1) and run formatting (lang="en-US").
$lastPara = $d.Paragraphs.Last
[void]$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newRange = $newPara.Range
$newRange.Collapse(1)

# 3. Fill the new paragraph with its final content: three runs of text, the
#    "Renderer.cs" run wrapped in spellStart/spellEnd proofErr markers, and
#    the relocated bookmark sitting between the second and third run.
$paraXml = (
    '<w:p>' +
        '<w:pPr>' +
            '<w:pStyle w:val="a3"/>' +
            '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' +
        '</w:pPr>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr>' +
            '<w:t xml:space="preserve">When listing empty players list exception fixed. Add Game constructor. </w:t>' +
        '</w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Renderer.cs</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr>' +
            '<w:t xml:space="preserve"> left with default constructor (used parameter injection </w:t>' +
        '</w:r>' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
        '<w:bookmarkEnd w:id="0"/>' +
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr>' +
            '<w:t xml:space="preserve">instead) </w:t>' +
        '</w:r>' +
    '</w:p>'
)
$packageXml = (
    '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
            '<pkg:xmlData>' +
                '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                    '<w:body>' + $paraXml + '</w:body>' +
                '</w:document>' +
            '</pkg:xmlData>' +
        '</pkg:part>' +
    '</pkg:package>'
)
[void]$newRange.InsertXML($packageXml)

# 4. Inserting paragraph XML at the very end of the document body leaves a
#    spare empty paragraph behind it; merge that leftover paragraph mark away
#    so the text we just inserted is once again the last paragraph of the
#    body (directly followed by the sectPr), matching the original layout.
$paraCount = $d.Paragraphs.Count
if ($paraCount -gt 0 -and $d.Paragraphs.Item($paraCount).Range.Text -eq "") {
    $secondLast = $d.Paragraphs.Item($paraCount - 1)
    $mergeRange = $d.Range($secondLast.Range.End - 1, $d.Content.End)
    [void]$mergeRange.Delete()
}
